# PROS-11827 - additional "Other Main Shelf" task added to the SOS and
# Availability KPI "Scene type / Tasks" list, plus the window/selection
# state and a stale external-workbook reference index that moved along
# with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Content edit -----------------------------------------------------
# Row 2 ("Facings SOS") and row 3 ("Availability") both list the same
# set of Main Shelf tasks in column C; append "Other Main Shelf" to both.
$newTasks = "Pain Main Shelf, Oral Main Shelf, Respiratory Main Shelf, NRT Main Shelf, Other Main Shelf"
$ws.Range("C2").Value = $newTasks
$ws.Range("C3").Value = $newTasks

# --- Fix stale external reference index --------------------------------
# The workbook only has one external reference left, so the
# Validation_List defined name should point at [1] (not the old [2]).
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Validation_List") {
        $n.RefersTo = "=[1]Set_up!`$A`$90:`$A`$124"
    }
}

# --- Window / selection state -------------------------------------------
# Scroll the frozen sheet so the visible pane starts at column E instead
# of column I, and leave the active selection on K19.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 2
$ws.Range("K19").Select()
